# Add data for 2025-07-05
# Updates crime-count cells (mostly the partial-year 2025 column "L", plus a couple
# of minor revisions to 2015 "B" / 2024 "K" columns) across the Citywide Totals sheet,
# the By Neighborhood rollup sheet, and the individual neighborhood sheets that had
# new/updated incidents on that date.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 3352
$ws.Range("L3").Value = 3466
$ws.Range("B4").Value = 1713
$ws.Range("K4").Value = 1768
$ws.Range("L4").Value = 869
$ws.Range("L5").Value = 205
$ws.Range("L6").Value = 3053
$ws.Range("B7").Value = 23345
$ws.Range("K7").Value = 27558
$ws.Range("L7").Value = 10945

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 94
$ws.Range("L4").Value = 43
$ws.Range("L5").Value = 42
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 366
$ws.Range("L8").Value = 699
$ws.Range("L9").Value = 70
$ws.Range("L10").Value = 69
$ws.Range("L11").Value = 185
$ws.Range("L14").Value = 57
$ws.Range("L15").Value = 82
$ws.Range("L18").Value = 86
$ws.Range("L19").Value = 306
$ws.Range("L20").Value = 280
$ws.Range("L23").Value = 119
$ws.Range("L24").Value = 27
$ws.Range("L26").Value = 11
$ws.Range("K27").Value = 261
$ws.Range("L29").Value = 591
$ws.Range("L31").Value = 103
$ws.Range("L33").Value = 506
$ws.Range("L36").Value = 146
$ws.Range("L37").Value = 391
$ws.Range("L42").Value = 349
$ws.Range("L43").Value = 82
$ws.Range("K51").Value = 357
$ws.Range("L52").Value = 215
$ws.Range("L53").Value = 126
$ws.Range("L54").Value = 229
$ws.Range("L55").Value = 105
$ws.Range("B63").Value = 417
$ws.Range("K63").Value = 161
$ws.Range("L63").Value = 39
$ws.Range("L65").Value = 211
$ws.Range("L69").Value = 29
$ws.Range("L73").Value = 93
$ws.Range("L76").Value = 157
$ws.Range("L79").Value = 279
$ws.Range("L83").Value = 257
$ws.Range("L84").Value = 109
$ws.Range("L85").Value = 554
$ws.Range("L88").Value = 128
$ws.Range("L89").Value = 150
$ws.Range("L91").Value = 155
$ws.Range("L93").Value = 60
$ws.Range("L94").Value = 133
$ws.Range("L95").Value = 154
$ws.Range("L96").Value = 107
$ws.Range("L99").Value = 182
$ws.Range("B101").Value = 23345
$ws.Range("K101").Value = 27558
$ws.Range("L101").Value = 10945

# Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 57

# West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 107

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 118
$ws.Range("L3").Value = 111
$ws.Range("L6").Value = 104
$ws.Range("L7").Value = 366

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 185

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 45
$ws.Range("L3").Value = 44
$ws.Range("L4").Value = 25
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 150

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 159
$ws.Range("L3").Value = 225
$ws.Range("L4").Value = 41
$ws.Range("L6").Value = 118
$ws.Range("L7").Value = 554

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 65
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 215

# Norwood Park
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 29

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 39
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 126

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 206
$ws.Range("L3").Value = 229
$ws.Range("L4").Value = 47
$ws.Range("L6").Value = 195
$ws.Range("L7").Value = 699

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 102
$ws.Range("L7").Value = 257

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 155
$ws.Range("L4").Value = 27
$ws.Range("L7").Value = 506

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 154

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 117
$ws.Range("L3").Value = 118
$ws.Range("L7").Value = 391

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 64
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 211

# Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 182

# Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 33
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 103

# South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 109

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L2").Value = 49
$ws.Range("L3").Value = 53
$ws.Range("L6").Value = 109
$ws.Range("L7").Value = 229

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 229
$ws.Range("L7").Value = 591

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 93
$ws.Range("L7").Value = 306

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 72
$ws.Range("L7").Value = 157

# Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L2").Value = 36
$ws.Range("L7").Value = 87

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 104
$ws.Range("L3").Value = 112
$ws.Range("L6").Value = 96
$ws.Range("L7").Value = 349

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 30
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 69

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 105

# Dunning
$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 27

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 119

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 155

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 101
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 279

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 92
$ws.Range("L4").Value = 26
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 280

# Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L2").Value = 33
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 86

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 146

# West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L2").Value = 21
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 60

# West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 133

# Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 82

# East Village
$ws = $wb.Worksheets.Item('East Village')
$ws.Range("L2").Value = 2
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 11

# Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 70

# Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 93
$ws.Range("L4").Value = 8

# Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L7").Value = 94

# United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L2").Value = 35
$ws.Range("L3").Value = 46
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 128

# Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 42

# Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K4").Value = 33
$ws.Range("K7").Value = 261

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K4").Value = 41
$ws.Range("K7").Value = 357

# Hyde Park
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 82

# Archer Heights
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 43

